# Apply the "Update dokumentace (export časopisu)" changes.
$d = $word.ActiveDocument

# 1) Version bump v0.2 -> v0.3
$d.Content.Find.Execute("v0.2", $true, $false, $false, $false, $false, `
    $true, 1, $false, "v0.3", 2) | Out-Null

# 2) "Čeká" -> "Čeká," in the "a 3) Redaktor ..." bullet, and drop the
#    surrounding gramStart/gramEnd proof-error markers by re-writing the
#    whole sentence (anchored uniquely by its preceding text).
$anchor = $d.Content
$anchor.Find.Execute("a 3) Redaktor může použít pouze základní funkce. Čeká než odešlou oba recenzenti posudky, aby jej mohl odeslat autorovi.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "a 3) Redaktor může použít pouze základní funkce. Čeká, než odešlou oba recenzenti posudky, aby jej mohl odeslat autorovi.", 2) | Out-Null

# 3) Insert a new sentence into the "Příspěvek je přijat k vydání" bullet.
$d.Content.Find.Execute("časopisu. Tuto volbu lze opětovně zrušit tlačítkem", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "časopisu. A může být exportován s ostatními vydanými články v časopisu v menu správa časopisu. Tuto volbu lze opětovně zrušit tlačítkem", 2) | Out-Null

# 4) Insert " a jeho článků" into the "U jednotlivých časopisů ..." sentence.
$d.Content.Find.Execute("zveřejnění časopisu veřejnosti na hlavní stránce.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "zveřejnění časopisu a jeho článků veřejnosti na hlavní stránce.", 2) | Out-Null

# 5) Rework the "Exportovat" sentence: describe the zip archive and drop the
#    "(Aktuálně neimplementováno)" note.
$d.Content.Find.Execute("časopisu a stáhne je. (Aktuálně neimplementováno)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "časopisu do komprimovaného souboru archivu (zip) a stáhne je.", 2) | Out-Null
